$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.819.18'
$ws.Range("E2").Value = '  -6.07%  '

$ws.Range("D3").Value = '2.540.49'
$ws.Range("E3").Value = '  -5.11%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '299.16'
$ws.Range("E5").Value = '  -3.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.78'
$ws.Range("E6").Value = '  -5.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.574'
$ws.Range("E7").Value = '  -4.21%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.548'
$ws.Range("E9").Value = '  -6.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.20'
$ws.Range("E10").Value = '  -6.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0805'
$ws.Range("E11").Value = '  -5.34%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.70'
$ws.Range("E12").Value = '  -5.52%  '

$ws.Range("E13").Value = '  +6.93%  '

$ws.Range("D14").Value = '2.925.66'
$ws.Range("E14").Value = '  -5.50%  '

$ws.Range("D15").Value = '2.588.38'
$ws.Range("E15").Value = '  -3.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.877'
$ws.Range("E16").Value = '  -6.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.18'
$ws.Range("E17").Value = '  -6.69%  '

$ws.Range("D18").Value = '42.778.90'
$ws.Range("E18").Value = '  -6.54%  '

$ws.Range("D19").Value = '0.0₃0980'
$ws.Range("E19").Value = '  -4.20%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.64'
$ws.Range("E20").Value = '  -1.89%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.53'
$ws.Range("E21").Value = '  -4.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.75'
$ws.Range("E22").Value = '  -4.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.98'
$ws.Range("E23").Value = '  -10.16%  '

$ws.Range("E24").Value = '  -5.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.13'
$ws.Range("E25").Value = '  -5.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.07'
$ws.Range("E26").Value = '  -5.82%  '

$ws.Range("E27").Value = '  +0.35%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.18'
$ws.Range("E28").Value = '  -4.53%  '

$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.66'
$ws.Range("E29").Value = '  -4.93%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.11'
$ws.Range("E30").Value = '  -4.89%  '

$ws.Range("E31").Value = '  -2.76%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '152.70'
$ws.Range("E32").Value = '  -1.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.76'
$ws.Range("E33").Value = '  -2.52%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.36'
$ws.Range("E34").Value = '  -10.62%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.14'
$ws.Range("E35").Value = '  -9.22%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0793'
$ws.Range("E36").Value = '  -6.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.114'
$ws.Range("E37").Value = '  -6.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '17.45'
$ws.Range("E38").Value = '  +7.07%  '

$ws.Range("E39").Value = '  -4.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.96'
$ws.Range("E40").Value = '  -11.18%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.41'
$ws.Range("E41").Value = '  -6.69%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0309'
$ws.Range("E42").Value = '  -5.90%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.83'
$ws.Range("E43").Value = '  -4.51%  '

$ws.Range("B44").Value = 'Maker'
$ws.Range("C44").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D44").Value = '2.086.57'
$ws.Range("E44").Value = '  -1.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.92'
$ws.Range("E45").Value = '  +20.42%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.997'
$ws.Range("E46").Value = '  -0.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.10'
$ws.Range("E47").Value = '  -2.60%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '84.22'
$ws.Range("E48").Value = '  -10.32%  '

$ws.Range("D49").Value = '2.785.57'
$ws.Range("E49").Value = '  -5.26%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.70'
$ws.Range("E50").Value = '  -6.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.67'
$ws.Range("E51").Value = '  -3.77%  '
